# feat: add 2022-Q4 data
#
# Current workbook layout:
#   Sheet1 "总计"    - summary sheet, one row per quarter snapshot
#   Sheet2 "2022-Q2" - per-fund detail for the 2022-Q2 snapshot
#
# Target workbook layout:
#   Sheet1 "总计"    - gains a new row for the 2022-Q4 snapshot (inserted
#                      ahead of the existing 2022-Q2 row)
#   Sheet2 "2022-Q4" - (re-uses the old "2022-Q2" sheet/tab) now holds the
#                      brand-new per-fund detail for 2022-Q4
#   Sheet3 "2022-Q2" - brand-new tab, an exact copy of the original
#                      "2022-Q2" detail sheet, appended at the end

$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell while forcing text storage (t="inlineStr"/
# shared-string) even when the text looks numeric (e.g. "0.44"), matching
# how these "numbers-as-text" columns are stored in the source data. The
# temporary "@" number format is cleared again right after via Style =
# "Normal" so no stray formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" detail sheet to the end of the
#    workbook before touching it, so its current content + formatting is
#    preserved verbatim on the new tab.
# ---------------------------------------------------------------------
$q2Sheet.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$q2Copy = $wb.Worksheets.Item($wb.Worksheets.Count)

# Free up the "2022-Q2" name on the original tab before renaming, then
# rename the copy back to it.
$q2Sheet.Name = "2022-Q4"
$q2Copy.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 2. "总计" (summary) sheet: relabel the existing data row as 2022-Q4 and
#    append a new row restating the old 2022-Q2 figures below it.
# ---------------------------------------------------------------------
$totalSheet.Range("B2").Value = "2022-Q4"

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

# ---------------------------------------------------------------------
# 3. Overwrite the "2022-Q4" detail sheet (former "2022-Q2" tab, now at
#    position 2) with the new Q4 fund data.
# ---------------------------------------------------------------------

# Header row formatting: match the "总计" sheet's header style.
$totalSheet.Range("B1").Copy()
$q2Sheet.Range("B1:H1").PasteSpecial(-4122)

$q2Sheet.Range("B1").Value = "基金代码"
$q2Sheet.Range("C1").Value = "基金名称"
$q2Sheet.Range("D1").Value = "基金规模"
$q2Sheet.Range("E1").Value = "股票总仓位"
$q2Sheet.Range("F1").Value = "仓位占比"
$q2Sheet.Range("G1").Value = "持有市值(亿元)"
$q2Sheet.Range("H1").Value = "仓位排名"

# Row 2 - style for A2 matches the "总计" sheet's A2 style.
$totalSheet.Range("A2").Copy()
$q2Sheet.Range("A2").PasteSpecial(-4122)
$q2Sheet.Range("A2").Value = 0
Set-TextValue $q2Sheet.Range("B2") "519981"
Set-TextValue $q2Sheet.Range("C2") "长信美国标准普尔100等权重指数增强（QDII）人民币"
Set-TextValue $q2Sheet.Range("D2") "0.44"
Set-TextValue $q2Sheet.Range("E2") "82.94"
Set-TextValue $q2Sheet.Range("F2") "0.87"
Set-TextValue $q2Sheet.Range("G2") "0.0038"
$q2Sheet.Range("H2").Value = 5

# Row 3
$totalSheet.Range("A2").Copy()
$q2Sheet.Range("A3").PasteSpecial(-4122)
$q2Sheet.Range("A3").Value = 1
Set-TextValue $q2Sheet.Range("B3") "011706"
Set-TextValue $q2Sheet.Range("C3") "长信美国标准普尔100等权重指数增强（QDII）美元"
Set-TextValue $q2Sheet.Range("D3") "0.44"
Set-TextValue $q2Sheet.Range("E3") "82.94"
Set-TextValue $q2Sheet.Range("F3") "0.87"
Set-TextValue $q2Sheet.Range("G3") "0.0038"
$q2Sheet.Range("H3").Value = 5

# Page margins for the new Q4 sheet differ from the old Q2 sheet's.
# PageSetup margins are expressed in points (72pt = 1 inch); the XML stores
# inches, so multiply the desired inch values by 72.
$q2Sheet.PageSetup.LeftMargin = 0.75 * 72
$q2Sheet.PageSetup.RightMargin = 0.75 * 72
$q2Sheet.PageSetup.TopMargin = 1 * 72
$q2Sheet.PageSetup.BottomMargin = 1 * 72
$q2Sheet.PageSetup.HeaderMargin = 0.5 * 72
$q2Sheet.PageSetup.FooterMargin = 0.5 * 72

$totalSheet.Activate()
